$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note text in A1 ---
$hoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.91 = 6953.26 pesos`n✅ 6953.26 pesos = 1.9 = 950.1 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$hoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate figures ---
$tasas = $wb.Worksheets.Item("tasas")

$tasas.Range("N10").Value = 522.776
$tasas.Range("O10").Value = 3635

$tasas.Range("N12").Value = 3651.99
$tasas.Range("O12").Value = 499.01
